$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: TestCase Number 9
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Validate_EmailSentToUser_WhileForgotPassword"
$ws.Cells.Item(10, 3).Value = "Functional"
$ws.Cells.Item(10, 4).Value = "Verify that email has been sent to user when user name is valid and existing while reset password."

# Row 11: TestCase Number 10
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Validate_SuccessfulLogin"
$ws.Cells.Item(11, 3).Value = "Functional"
$ws.Cells.Item(11, 4).Value = "Verif successful login when email and password is valid and existing in database of kirana bazaar."

# Match the style (wrap text) of column D used by the other description cells
$ws.Range("D10").WrapText = $true
$ws.Range("D11").WrapText = $true

# Match row height used by the other wrapped-text data rows (30pt)
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# Update selection to match the committed state
$ws.Range("D11").Select()
